# Auto-generated edit script applying numeric updates to Kujata_Profits sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (42 cell changes) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 360.66666
$ws.Range("I5").Value = 305.75
$ws.Range("J5").Value = 470.5
$ws.Range("K5").Value = 305.75
$ws.Range("L5").Value = 470.5
$ws.Range("M5").Value = -190.75
$ws.Range("N5").Value = -700.5
$ws.Range("H107").Value = 4209.5
$ws.Range("I107").Value = 3176.6428
$ws.Range("J107").Value = 7824.5
$ws.Range("K107").Value = 3176.6428
$ws.Range("L107").Value = 7824.5
$ws.Range("M107").Value = -1256.6428
$ws.Range("N107").Value = -11664.5
$ws.Range("H132").Value = 7253809.5
$ws.Range("I132").Value = 8338271
$ws.Range("J132").Value = 24067.666
$ws.Range("K132").Value = 25014813
$ws.Range("L132").Value = 72202.99800000001
$ws.Range("M132").Value = -25012283
$ws.Range("N132").Value = -77262.99800000001
$ws.Range("H135").Value = 40000896
$ws.Range("I135").Value = 638.05884
$ws.Range("J135").Value = 125001440
$ws.Range("K135").Value = 5742.52956
$ws.Range("L135").Value = 1125012960
$ws.Range("M135").Value = -3207.52956
$ws.Range("N135").Value = -1125018030
$ws.Range("H137").Value = 2089.0188
$ws.Range("I137").Value = 1689.7693
$ws.Range("J137").Value = 2473.4814
$ws.Range("K137").Value = 5069.3079
$ws.Range("L137").Value = 7420.4442
$ws.Range("M137").Value = -2519.3079
$ws.Range("N137").Value = -12520.4442
$ws.Range("H138").Value = 2205.6
$ws.Range("I138").Value = 978.06665
$ws.Range("J138").Value = 2422.2236
$ws.Range("K138").Value = 2934.19995
$ws.Range("L138").Value = 7266.6708
$ws.Range("M138").Value = 2205.80005
$ws.Range("N138").Value = -17546.6708

# ---- Sheet: ARM (32 cell changes) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 1500
$ws.Range("I16").Value = 1500
$ws.Range("K16").Value = 1500
$ws.Range("M16").Value = -1213
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 74005
$ws.Range("I23").Value = 67504.5
$ws.Range("K23").Value = 67504.5
$ws.Range("M23").Value = -67245.5
$ws.Range("H26").Value = 77.5
$ws.Range("I26").Value = 77.5
$ws.Range("K26").Value = 77.5
$ws.Range("M26").Value = 252.5
$ws.Range("H36").Value = 1200
$ws.Range("I36").Value = 1200
$ws.Range("K36").Value = 1200
$ws.Range("M36").Value = -854
$ws.Range("H88").Value = 2349.875
$ws.Range("I88").Value = 2749.5
$ws.Range("K88").Value = 2749.5
$ws.Range("M88").Value = -2343.5
$ws.Range("H91").Value = 2349.875
$ws.Range("I91").Value = 2749.5
$ws.Range("K91").Value = 2749.5
$ws.Range("M91").Value = -1345.5
$ws.Range("H114").Value = 40666.332
$ws.Range("J114").Value = 40666.332
$ws.Range("L114").Value = 40666.332
$ws.Range("N114").Value = -49344.332

# ---- Sheet: BSM (29 cell changes) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3736.55
$ws.Range("I86").Value = 3808.9285
$ws.Range("J86").Value = 3567.6667
$ws.Range("K86").Value = 3808.9285
$ws.Range("L86").Value = 3567.6667
$ws.Range("M86").Value = -2685.9285
$ws.Range("N86").Value = -5813.6667
$ws.Range("H89").Value = 3736.55
$ws.Range("I89").Value = 3808.9285
$ws.Range("J89").Value = 3567.6667
$ws.Range("K89").Value = 19044.6425
$ws.Range("L89").Value = 17838.3335
$ws.Range("M89").Value = -13428.6425
$ws.Range("N89").Value = -29070.3335
$ws.Range("H107").Value = 1434.7778
$ws.Range("I107").Value = 1171.4286
$ws.Range("J107").Value = 2356.5
$ws.Range("K107").Value = 1171.4286
$ws.Range("L107").Value = 2356.5
$ws.Range("M107").Value = 748.5714
$ws.Range("N107").Value = -6196.5
$ws.Range("H110").Value = 19500
$ws.Range("J110").Value = 19500
$ws.Range("L110").Value = 19500
$ws.Range("N110").Value = -27680
$ws.Range("H140").Value = 27949.955
$ws.Range("J140").Value = 27949.955
$ws.Range("L140").Value = 27949.955
$ws.Range("N140").Value = -38309.955

# ---- Sheet: CRP (30 cell changes) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1332.7255
$ws.Range("I31").Value = 1290.1459
$ws.Range("K31").Value = 1290.1459
$ws.Range("M31").Value = -995.1459
$ws.Range("H34").Value = 1332.7255
$ws.Range("I34").Value = 1290.1459
$ws.Range("K34").Value = 1290.1459
$ws.Range("M34").Value = -1088.1459
$ws.Range("H58").Value = 9186.5625
$ws.Range("I58").Value = 1996.6666
$ws.Range("J58").Value = 10845.77
$ws.Range("K58").Value = 1996.6666
$ws.Range("L58").Value = 10845.77
$ws.Range("M58").Value = -1793.6666
$ws.Range("N58").Value = -11251.77
$ws.Range("H132").Value = 2301
$ws.Range("I132").Value = 2184.3333
$ws.Range("K132").Value = 6552.999899999999
$ws.Range("M132").Value = -4022.999899999999
$ws.Range("H135").Value = 36640.715
$ws.Range("J135").Value = 36640.715
$ws.Range("L135").Value = 36640.715
$ws.Range("N135").Value = -46780.715
$ws.Range("H136").Value = 9186.5625
$ws.Range("I136").Value = 1996.6666
$ws.Range("J136").Value = 10845.77
$ws.Range("K136").Value = 5989.9998
$ws.Range("L136").Value = 32537.31
$ws.Range("M136").Value = -3439.9998
$ws.Range("N136").Value = -37637.31

# ---- Sheet: CUL (37 cell changes) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4787227
$ws.Range("J4").Value = 5972357
$ws.Range("L4").Value = 17917071
$ws.Range("N4").Value = -17917295
$ws.Range("H14").Value = 232.125
$ws.Range("I14").Value = 232.125
$ws.Range("K14").Value = 696.375
$ws.Range("M14").Value = -523.375
$ws.Range("H86").Value = 836.6667
$ws.Range("I86").Value = 671.4286
$ws.Range("J86").Value = 981.25
$ws.Range("K86").Value = 2014.2858
$ws.Range("L86").Value = 2943.75
$ws.Range("N86").Value = -5315.75
$ws.Range("M86").Value = -828.2857999999999
$ws.Range("H87").Value = 1287.6
$ws.Range("J87").Value = 1600
$ws.Range("L87").Value = 4800
$ws.Range("N87").Value = -7296
$ws.Range("H89").Value = 836.6667
$ws.Range("I89").Value = 671.4286
$ws.Range("J89").Value = 981.25
$ws.Range("K89").Value = 6042.8574
$ws.Range("L89").Value = 8831.25
$ws.Range("N89").Value = -20687.25
$ws.Range("M89").Value = -114.8573999999999
$ws.Range("H90").Value = 1287.6
$ws.Range("J90").Value = 1600
$ws.Range("L90").Value = 14400
$ws.Range("N90").Value = -26880
$ws.Range("H131").Value = 31297902
$ws.Range("I131").Value = 90909750
$ws.Range("J131").Value = 72647.664
$ws.Range("K131").Value = 272729250
$ws.Range("L131").Value = 217942.992
$ws.Range("M131").Value = -272724210
$ws.Range("N131").Value = -228022.992

# ---- Sheet: GSM (33 cell changes) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10333.5
$ws.Range("I5").Value = 1000.5
$ws.Range("J5").Value = 15000
$ws.Range("K5").Value = 1000.5
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = -888.5
$ws.Range("N5").Value = -15224
$ws.Range("H80").Value = 3566.6667
$ws.Range("J80").Value = 5775
$ws.Range("L80").Value = 5775
$ws.Range("N80").Value = -7771
$ws.Range("H83").Value = 3566.6667
$ws.Range("J83").Value = 5775
$ws.Range("L83").Value = 28875
$ws.Range("N83").Value = -38859
$ws.Range("H126").Value = 2225.2727
$ws.Range("I126").Value = 1839.8572
$ws.Range("J126").Value = 2899.75
$ws.Range("K126").Value = 5519.571599999999
$ws.Range("L126").Value = 8699.25
$ws.Range("M126").Value = -3049.571599999999
$ws.Range("N126").Value = -13639.25
$ws.Range("H132").Value = 3696.1052
$ws.Range("I132").Value = 3654.6365
$ws.Range("J132").Value = 3753.125
$ws.Range("K132").Value = 10963.9095
$ws.Range("L132").Value = 11259.375
$ws.Range("M132").Value = -8433.9095
$ws.Range("N132").Value = -16319.375
$ws.Range("H139").Value = 36968
$ws.Range("J139").Value = 36968
$ws.Range("L139").Value = 36968
$ws.Range("N139").Value = -47248

# ---- Sheet: LTW (30 cell changes) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 816
$ws.Range("I22").Value = 686.2222
$ws.Range("J22").Value = 1400
$ws.Range("K22").Value = 686.2222
$ws.Range("L22").Value = 1400
$ws.Range("M22").Value = -391.2222
$ws.Range("N22").Value = -1990
$ws.Range("H27").Value = 816
$ws.Range("I27").Value = 686.2222
$ws.Range("J27").Value = 1400
$ws.Range("K27").Value = 686.2222
$ws.Range("L27").Value = 1400
$ws.Range("M27").Value = -579.2222
$ws.Range("N27").Value = -1614
$ws.Range("H82").Value = 2211.5715
$ws.Range("I82").Value = 2246.8333
$ws.Range("K82").Value = 2246.8333
$ws.Range("M82").Value = -1885.8333
$ws.Range("H85").Value = 2211.5715
$ws.Range("I85").Value = 2246.8333
$ws.Range("K85").Value = 2246.8333
$ws.Range("M85").Value = -998.8332999999998
$ws.Range("H93").Value = 1018.6667
$ws.Range("I93").Value = 1018.6667
$ws.Range("K93").Value = 1018.6667
$ws.Range("M93").Value = 229.3333
$ws.Range("H136").Value = 2028.625
$ws.Range("J136").Value = 2321.5
$ws.Range("L136").Value = 6964.5
$ws.Range("N136").Value = -12064.5

# ---- Sheet: WVR (29 cell changes) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 12000
$ws.Range("J94").Value = 12000
$ws.Range("L94").Value = 12000
$ws.Range("N94").Value = -13802
$ws.Range("H100").Value = 450
$ws.Range("I100").Value = 450
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 900
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -359
$ws.Range("N100").ClearContents()
$ws.Range("H119").Value = 10786.4
$ws.Range("J119").Value = 10824.5
$ws.Range("L119").Value = 10824.5
$ws.Range("N119").Value = -20500.5
$ws.Range("H132").Value = 2459.4243
$ws.Range("I132").Value = 2178.0344
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 6534.1032
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -4004.1032
$ws.Range("N132").Value = -18558.5
$ws.Range("H136").Value = 1736.3158
$ws.Range("I136").Value = 1430.0769
$ws.Range("J136").Value = 2399.8333
$ws.Range("K136").Value = 4290.2307
$ws.Range("L136").Value = 7199.499899999999
$ws.Range("M136").Value = -1740.2307
$ws.Range("N136").Value = -12299.4999

Write-Host "Applied all cell updates."